# add save column in s_vals sheets
# The sheet has headers in row 1 (B1:G1) and data in rows 2-3 (A2:G3).
# "sum" is the last header column (G). We append a new "Save" column in H,
# filled with 1 for every data row, copying the header's formatting from
# the existing "sum" header cell so the new column matches the sheet style.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Find the last used row/column on the sheet so this generalises to any
# s_vals-style sheet (header row 1, data starting row 2).
$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count
$lastCol = $usedRange.Columns.Count

$headerCell = $ws.Cells.Item(1, $lastCol)
$newCol = $lastCol + 1
$newHeaderCell = $ws.Cells.Item(1, $newCol)

# Copy formatting (style/border/font/alignment) from the existing header
# cell onto the new header cell, then set its text.
$headerCell.Copy($newHeaderCell)
$newHeaderCell.Value = "Save"

# Fill the new column's data rows with 1.
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, $newCol).Value = 1
}
